$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.781.62"
$ws.Range("E2").Value = "  +3.00%  "

$ws.Range("D3").Value = "1.866.28"
$ws.Range("E3").Value = "  +2.79%  "

$ws.Range("D4").Value = "'1.040"
$ws.Range("E4").Value = "  +3.53%  "

$ws.Range("D5").Value = "'324.55"
$ws.Range("E5").Value = "  +3.98%  "

$ws.Range("E6").Value = "  +3.09%  "

$ws.Range("D7").Value = "'0.4427"
$ws.Range("E7").Value = "  +3.03%  "

$ws.Range("D8").Value = "'0.3804"
$ws.Range("E8").Value = "  +3.75%  "

$ws.Range("D9").Value = "'0.07474"
$ws.Range("E9").Value = "  +3.30%  "

$ws.Range("D10").Value = "'0.8858"
$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D11").Value = "'21.73"
$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").Value = "1.885.27"
$ws.Range("E12").Value = "  -13.20%  "

$ws.Range("D13").Value = "'5.565"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("D14").Value = "'6.765"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").Value = "'0.07242"
$ws.Range("E15").Value = "  +4.20%  "

$ws.Range("D16").Value = "'83.90"
$ws.Range("E16").Value = "  +3.39%  "

$ws.Range("D17").Value = "'1.042"
$ws.Range("E17").Value = "  +3.57%  "

$ws.Range("D18").Value = "'0.000009159"
$ws.Range("E18").Value = "  +3.04%  "

$ws.Range("D19").Value = "'1.036"
$ws.Range("E19").Value = "  +3.04%  "

$ws.Range("E20").Value = "  +2.65%  "

$ws.Range("D21").Value = "27.799.43"
$ws.Range("E21").Value = "  +2.90%  "

$ws.Range("D22").Value = "'5.327"
$ws.Range("E22").Value = "  +3.00%  "

$ws.Range("E23").Value = "  +3.48%  "

$ws.Range("D24").Value = "'1.980"
$ws.Range("E24").Value = "  +5.06%  "

$ws.Range("D25").Value = "'158.92"
$ws.Range("E25").Value = "  +3.26%  "

$ws.Range("D26").Value = "'18.89"
$ws.Range("E26").Value = "  +3.04%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'1.989"
$ws.Range("E27").Value = "  +4.85%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'5.332"
$ws.Range("E28").Value = "  +2.08%  "

$ws.Range("D29").Value = "'117.90"
$ws.Range("E29").Value = "  +2.86%  "

$ws.Range("D30").Value = "'0.09083"
$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("D31").Value = "'3.119"
$ws.Range("E31").Value = "  +11.18%  "

$ws.Range("D32").Value = "'0.7786"
$ws.Range("E32").Value = "  +4.39%  "

$ws.Range("D33").Value = "'1.215"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D34").Value = "'4.578"
$ws.Range("E34").Value = "  +3.69%  "

$ws.Range("E35").Value = "  +3.16%  "

$ws.Range("D36").Value = "'1.159"
$ws.Range("E36").Value = "  +2.47%  "

$ws.Range("D37").Value = "'0.01996"
$ws.Range("E37").Value = "  +3.89%  "

$ws.Range("E38").Value = "  +2.85%  "

$ws.Range("D39").Value = "'2.869"
$ws.Range("E39").Value = "  +4.77%  "

$ws.Range("D40").Value = "'0.5207"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("D41").Value = "'0.1695"

$ws.Range("D42").Value = "'6.897"
$ws.Range("E42").Value = "  +7.03%  "

$ws.Range("D43").Value = "'8.694"
$ws.Range("E43").Value = "  +4.72%  "

$ws.Range("E44").Value = "  +3.73%  "

$ws.Range("D45").Value = "'109.76"
$ws.Range("E45").Value = "  +2.88%  "

$ws.Range("D46").Value = "'1.723"
$ws.Range("E46").Value = "  +4.96%  "

$ws.Range("D47").Value = "'0.4716"
$ws.Range("E47").Value = "  +2.94%  "

$ws.Range("E48").Value = "  +4.32%  "

$ws.Range("D49").Value = "'1.918"
$ws.Range("E49").Value = "  +4.32%  "

$ws.Range("D50").Value = "'39.94"
$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("D51").Value = "'64.66"
$ws.Range("E51").Value = "  +2.70%  "
